$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "299.94"
Set-TextValue "E2" "-0.27%"

# Row 3 - OKB
Set-TextValue "D3" "31.80"
Set-TextValue "E3" "1.42%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.107"
Set-TextValue "E4" "-0.45%"

# Row 5 - Cronos
Set-TextValue "D5" "0.08207"
Set-TextValue "E5" "11.31%"

# Row 6 - FTXToken
Set-TextValue "D6" "2.582"
Set-TextValue "E6" "8.25%"

# Row 7 - KuCoinToken
Set-TextValue "D7" "7.846"
Set-TextValue "E7" "-1.47%"

# Row 8 - GateToken
Set-TextValue "D8" "3.842"
Set-TextValue "E8" "1.36%"

# Row 9 - MXToken
Set-TextValue "D9" "0.9296"
Set-TextValue "E9" "1.03%"

# Row 10 - WazirX
Set-TextValue "D10" "0.1763"
Set-TextValue "E10" "2.77%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-TextValue "D11" "0.07505"
Set-TextValue "E11" "-2.08%"

# Row 12 - MandalaExchangeToken
Set-TextValue "D12" "0.09009"
Set-TextValue "E12" "10.88%"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.03004"
Set-TextValue "E13" "-0.40%"

# Row 14 - BitMartToken (D unchanged)
Set-TextValue "E14" "0.99%"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001506"
Set-TextValue "E15" "0.76%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.005852"
Set-TextValue "E16" "-5.51%"

# Row 17 - LEO
Set-TextValue "D17" "3.613"
Set-TextValue "E17" "4.29%"

# Row 18 - BTSEToken
Set-TextValue "D18" "2.285"
Set-TextValue "E18" "2.64%"

# Row 19 - BitpandaEcosystemToken (D unchanged)
Set-TextValue "E19" "-1.16%"

# Row 21 - MCDex
Set-TextValue "D21" "3.905"
Set-TextValue "E21" "-16.09%"

# Row 22 - ZBToken (D unchanged)
Set-TextValue "E22" "7.00%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04624"
Set-TextValue "E23" "-0.45%"

# Row 24 - BitKan
Set-TextValue "D24" "0.001246"
Set-TextValue "E24" "1.76%"

# Row 25 - HotbitToken
Set-TextValue "D25" "0.004561"
Set-TextValue "E25" "1.68%"

# Row 26 - NitroEx
Set-TextValue "D26" "0.0001198"
Set-TextValue "E26" "-7.80%"

# Row 27 - UpBots
Set-TextValue "D27" "0.0003400"
Set-TextValue "E27" "81.67%"

# Row 39 - One
Set-TextValue "D39" "0.01778"
Set-TextValue "E39" "2.49%"

# Row 40 - IDEX
Set-TextValue "D40" "0.04586"
Set-TextValue "E40" "1.31%"

# Row 41 - KickToken
Set-TextValue "D41" "0.006904"
Set-TextValue "E41" "-4.40%"

# Row 42 - BKEXToken (E unchanged)
Set-TextValue "D42" "0.1380"

# Row 43 - CEJI
Set-TextValue "D43" "0.002207"
Set-TextValue "E43" "-0.11%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.01021"
Set-TextValue "E44" "-4.75%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00006182"
Set-TextValue "E45" "-1.43%"

# Row 46 - Kangarootoken
Set-TextValue "D46" "0.00000000749"
Set-TextValue "E46" "-0.17%"

# Row 47 - was CoinbaseStockToken, now BOLO
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D47" "0.7868"
Set-TextValue "E47" "-59.20%"

# Row 48 - was BOLO, now CoinbaseStockToken
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D48" "0.008375"
Set-TextValue "E48" "-16.27%"

# Row 49 - CryptobidCoin
Set-TextValue "D49" "0.00002097"
Set-TextValue "E49" "-0.17%"

# Row 50 - SpecialPowerGold
Set-TextValue "D50" "0.0001997"
Set-TextValue "E50" "-0.10%"
